$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.4370000000000001
$ws.Range("E2").Value = 0.1662
$ws.Range("G2").Value = 0.1234793187347932
$ws.Range("H2").Value = 0.1234793187347932
$ws.Range("I2").Value = 0.1058394160583942
$ws.Range("J2").Value = 0.09354965160592663
$ws.Range("K2").Value = 4.108
$ws.Range("L2").Value = 0.2498783454987835
$ws.Range("M2").Value = 1.9486
$ws.Range("N2").Value = 0.0360518038852914
$ws.Range("O2").Value = 0.4743427458617333
$ws.Range("P2").Value = 1.6926
$ws.Range("Q2").Value = 0.0313154486586494
$ws.Range("R2").Value = 0.4120253164556963
$ws.Range("S2").Value = 0.256
$ws.Range("T2").Value = 0.1313763727804578
$ws.Range("U2").Value = 57.64
$ws.Range("V2").Value = 1.066419981498612
$ws.Range("W2").Value = 0.148494983277592
$ws.Range("X2").Value = 0.019631896618046
$ws.Range("Y2").Value = 0.128863086659546
$ws.Range("Z2").Value = 6.277205040091617
$ws.Range("AB2").Value = 0.01967830629871981
$ws.Range("AC2").Value = -0.02617687185052062
$ws.Range("AD2").Value = 12.21
$ws.Range("AF2").Value = 12.21
$ws.Range("AG2").Value = -45.43
$ws.Range("AH2").Value = 0.1842740718382131
$ws.Range("AI2").Value = 0.2169509594882729
$ws.Range("AJ2").Value = -5.270301624129932
$ws.Range("AK2").Value = 33.4044117647059
$ws.Range("AL2").Value = 0.005
$ws.Range("AM2").Value = -0.045
$ws.Range("AN2").Value = 5.927184466019417
$ws.Range("AO2").Value = 348
$ws.Range("AP2").Value = -22.05339805825243
$ws.Range("AQ2").Value = -38.66666666666667
$ws.Range("D3").Value = 0.0548
$ws.Range("E3").Value = 0.0344
$ws.Range("K3").Value = 0.888
$ws.Range("L3").Value = 0.2187192118226601
$ws.Range("M3").Value = 0.336
$ws.Range("N3").Value = 0.03411167512690356
$ws.Range("O3").Value = 0.3783783783783784
$ws.Range("P3").Value = 0.08
$ws.Range("Q3").Value = 0.008121827411167513
$ws.Range("R3").Value = 0.09009009009009009
$ws.Range("S3").Value = 0.256
$ws.Range("T3").Value = 0.7619047619047619
$ws.Range("U3").Value = 1.62
$ws.Range("V3").Value = 0.1644670050761421
$ws.Range("W3").Value = 0.148494983277592
$ws.Range("X3").Value = 0.01960803972235912
$ws.Range("Y3").Value = 0.1288869435552328
$ws.Range("Z3").Value = 0.8144433299899696
$ws.Range("AB3").Value = 0.01962800848496591
$ws.Range("AC3").Value = -0.01962800848496591
$ws.Range("AD3").Value = 0.019
$ws.Range("AF3").Value = 0.019
$ws.Range("AG3").Value = -1.601
$ws.Range("AH3").Value = 0.001925220387070625
$ws.Range("AI3").Value = 0.002928032054245646
$ws.Range("AJ3").Value = -0.1940841314098679
$ws.Range("AK3").Value = -0.3288149517354694
$ws.Range("D4").Value = 0.526
$ws.Range("E4").Value = 0.298
$ws.Range("K4").Value = 1.67
$ws.Range("L4").Value = 0.2737704918032787
$ws.Range("M4").Value = 1.5006
$ws.Range("N4").Value = 0.06252500000000001
$ws.Range("O4").Value = 0.8985628742514972
$ws.Range("P4").Value = 1.5006
$ws.Range("Q4").Value = 0.06252500000000001
$ws.Range("R4").Value = 0.8985628742514972
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 4.22
$ws.Range("V4").Value = 0.1758333333333333
$ws.Range("W4").Value = 0.07167381974248926
$ws.Range("X4").Value = 0.02425934271993178
$ws.Range("Y4").Value = 0.04741447702255749
$ws.Range("Z4").Value = 0.2060810810810811
$ws.Range("AB4").Value = 0.02617687185052062
$ws.Range("AC4").Value = -0.02617687185052062
$ws.Range("AD4").Value = 12.1
$ws.Range("AF4").Value = 12.1
$ws.Range("AG4").Value = 7.88
$ws.Range("AH4").Value = 0.335180055401662
$ws.Range("AI4").Value = 0.3235294117647059
$ws.Range("AJ4").Value = 0.247176913425345
$ws.Range("AK4").Value = 0.2374924653405666
$ws.Range("AM4").Value = -0.005
$ws.Range("D5").Value = 0.4370000000000001
$ws.Range("G5").Value = 0.3232484076433121
$ws.Range("H5").Value = 0.3232484076433121
$ws.Range("I5").Value = 0.2770700636942675
$ws.Range("J5").Value = 0.238895966029724
$ws.Range("K5").Value = 1.55
$ws.Range("L5").Value = 0.2468152866242038
$ws.Range("M5").Value = 0.112
$ws.Range("N5").Value = 0.005544554455445545
$ws.Range("O5").Value = 0.07225806451612903
$ws.Range("P5").Value = 0.112
$ws.Range("Q5").Value = 0.005544554455445545
$ws.Range("R5").Value = 0.07225806451612903
$ws.Range("U5").Value = 51.8
$ws.Range("V5").Value = 2.564356435643564
$ws.Range("W5").Value = 0.1710816777041942
$ws.Range("X5").Value = 0.019631896618046
$ws.Range("Y5").Value = 0.1514497810861483
$ws.Range("Z5").Value = -0.1964587374084966
$ws.Range("AA5").Value = -0.04693319985818266
$ws.Range("AB5").Value = 0.01967830629871981
$ws.Range("AC5").Value = -0.06661150615690248
$ws.Range("AD5").Value = 0.091
$ws.Range("AF5").Value = 0.091
$ws.Range("AG5").Value = -51.709
$ws.Range("AH5").Value = 0.004484746932137401
$ws.Range("AI5").Value = 0.007344040029053345
$ws.Range("AJ5").Value = 1.641086673648799
$ws.Range("AK5").Value = 1.312111446623868
$ws.Range("AL5").Value = 0.005
$ws.Range("AM5").Value = -0.04
$ws.Range("AN5").Value = 0.04417475728155339
$ws.Range("AO5").Value = 348
$ws.Range("AP5").Value = -25.10145631067961
$ws.Range("AQ5").Value = -43.5
$ws.Range("E5").ClearContents()
